# VCEA_goals.xlsx — rename header columns to match naming conventions
# (drop "(Megawatts)" -> "_MW", drop " (%)" suffix from rps columns)
# and update the active selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "apco_rps"
$ws.Range("C1").Value = "dominion_rps "
$ws.Range("D1").Value = "apco_onshore_wind_and_solar_MW"
$ws.Range("E1").Value = "dominion_onshore_wind_and_solar_MW"
$ws.Range("F1").Value = "apco_storage_MW"
$ws.Range("G1").Value = "dominion_storage_MW"

# Move the selection/scroll state from D5 to C10 (and drop the
# topLeftCell="D1" frozen scroll position from the old view).
$ws.Activate()
$ws.Range("C10").Select()
